$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows (14-16) for the new RemitAcct2/3/4 fields, pushing the
# existing CreateDate/CreateEmpNo/LastUpdate/LastUpdateEmpNo rows down to 17-20.
$ws.Rows("14:16").Insert()

# Copy the formatting of the row above (row 13) down into the newly
# inserted rows so borders/number formats match the rest of the table.
$ws.Range("A13:E13").Copy($ws.Range("A14:E16"))

# SEQ numbers for the 3 new rows
$ws.Range("A14").Value = 6
$ws.Range("A15").Value = 7
$ws.Range("A16").Value = 8

# 欄位名稱 (field names) for the new rows
$ws.Range("B14").Value = "RemitAcct2"
$ws.Range("B15").Value = "RemitAcct3"
$ws.Range("B16").Value = "RemitAcct4"

# 中文名稱 (descriptions) for the new rows
$ws.Range("C14").Value = "調解匯款帳號"
$ws.Range("C15").Value = "更生匯款帳號"
$ws.Range("C16").Value = "清算匯款帳號"

# 形態 / 長度
$ws.Range("D14").Value = "VARCHAR2"
$ws.Range("D15").Value = "VARCHAR2"
$ws.Range("D16").Value = "VARCHAR2"
$ws.Range("E14").Value = 16
$ws.Range("E15").Value = 16
$ws.Range("E16").Value = 16

# Row 12: "RemitAcct" 的中文名稱 now reads "債協匯款帳號" instead of "匯款帳號"
$ws.Range("C12").Value = "債協匯款帳號"

# 備註說明 for the 3 new rows
$ws.Range("G14").Value = "20221122新增"
$ws.Range("G15").Value = "20221122新增"
$ws.Range("G16").Value = "20221122新增"

# Renumber the SEQ column for the rows that were pushed down
$ws.Range("A17").Value = 9
$ws.Range("A18").Value = 10
$ws.Range("A19").Value = 11
$ws.Range("A20").Value = 12

# Column G is a bit wider now to fit the new note text, and no longer "best fit"
$ws.Columns("G").ColumnWidth = 15.109375

# Selection moved
$ws.Range("H10").Select()
